# New weekly price observation for "Cebollín" at Terminal Hortofrutícola Agro
# Chillán is inserted as row 20 (sheet is sorted oldest price data last, with
# newest entries prepended near the top of the date-ordered block). Inserting
# a whole row shifts every existing row from 20-31 down to 21-32, which keeps
# all of their data intact while making room for the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("20").Insert()

# Fill in the newly inserted row 20 with the new observation's data.
$ws.Range("A20").Value = 7
$ws.Range("B20").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C20").Value = "Ñuble"
$ws.Range("D20").Value2 = 44839
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = 100112037
$ws.Range("G20").Value = "Cebollín"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 120
$ws.Range("K20").Value = 7500
$ws.Range("L20").Value = 8000
$ws.Range("M20").Value = 7750
$ws.Range("N20").Value = "$/docena de atados"
$ws.Range("O20").Value = "Provincia de Diguillín"
$ws.Range("P20").Value = 2583
$ws.Range("Q20").Value = 3
$ws.Range("R20").Value = "Hortaliza"
